$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.829.23'
$ws.Range("E2").Value = '  -0.52%  '

$ws.Range("D3").Value = '1.629.05'
$ws.Range("E3").Value = '  -0.49%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.14'
$ws.Range("E5").Value = '  +0.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5116'
$ws.Range("E6").Value = '  +0.61%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  +0.23%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06343'
$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.50'
$ws.Range("E10").Value = '  -0.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07783'
$ws.Range("E11").Value = '  +0.50%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.255'
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").Value = '1.633.63'
$ws.Range("E13").Value = '  -0.56%  '

$ws.Range("D14").Value = '1.849.44'
$ws.Range("E14").Value = '  -0.77%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5530'
$ws.Range("E15").Value = '  +1.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.65'
$ws.Range("E16").Value = '  -0.50%  '

$ws.Range("D17").Value = '0.0₅7511'
$ws.Range("E17").Value = '  -2.54%  '

$ws.Range("D18").Value = '25.806.59'

$ws.Range("E19").Value = '  -0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.431'
$ws.Range("E20").Value = '  +0.43%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '194.72'
$ws.Range("E21").Value = '  -1.99%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.790'
$ws.Range("E22").Value = '  -1.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.020'
$ws.Range("E23").Value = '  -0.33%  '

$ws.Range("E24").Value = '  -0.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.882'
$ws.Range("E25").Value = '  -0.40%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.47'
$ws.Range("E26").Value = '  +0.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1250'
$ws.Range("E27").Value = '  +4.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.57'
$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.722'
$ws.Range("E29").Value = '  -1.48%  '

$ws.Range("E30").Value = '  +0.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04870'
$ws.Range("E31").Value = '  -0.49%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.260'
$ws.Range("E32").Value = '  +0.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.174'
$ws.Range("E33").Value = '  +0.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.545'
$ws.Range("E34").Value = '  +1.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.358'
$ws.Range("E35").Value = '  -0.45%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8964'
$ws.Range("E36").Value = '  -1.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5527'
$ws.Range("E37").Value = '  +1.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.543'
$ws.Range("E38").Value = '  -1.71%  '

$ws.Range("D39").Value = '1.116.74'
$ws.Range("E39").Value = '  -0.87%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01551'
$ws.Range("E40").Value = '  -0.50%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9987'
$ws.Range("E41").Value = '  -0.19%  '

$ws.Range("E42").Value = '  +2.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7971'
$ws.Range("E43").Value = '  -1.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.46'
$ws.Range("E44").Value = '  -1.57%  '

$ws.Range("D45").Value = '1.775.50'
$ws.Range("E45").Value = '  +0.01%  '

$ws.Range("D46").Value = '0.0₈115'
$ws.Range("E46").Value = '  -7.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4425'
$ws.Range("E47").Value = '  -2.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9986'
$ws.Range("E48").Value = '  -0.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.66'
$ws.Range("E49").Value = '  -0.43%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05122'
$ws.Range("E50").Value = '  +0.20%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.589'
$ws.Range("E51").Value = '  +3.28%  '
